$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2940.875
$ws.Range("I40").Value = 2710
$ws.Range("J40").Value = 3079.4
$ws.Range("K40").Value = 2710
$ws.Range("L40").Value = 3079.4
$ws.Range("M40").Value = -2535
$ws.Range("N40").Value = -3429.4
$ws.Range("H53").Value = 2769
$ws.Range("I53").Value = 88.09999999999999
$ws.Range("J53").Value = 4831.231
$ws.Range("K53").Value = 88.09999999999999
$ws.Range("L53").Value = 4831.231
$ws.Range("M53").Value = 548.9
$ws.Range("N53").Value = -6105.231
$ws.Range("H55").Value = 353.22223
$ws.Range("I55").Value = 305
$ws.Range("J55").Value = 391.8
$ws.Range("K55").Value = 305
$ws.Range("L55").Value = 391.8
$ws.Range("M55").Value = -91
$ws.Range("N55").Value = -819.8
$ws.Range("H64").Value = 4378.5713
$ws.Range("I64").Value = 3860
$ws.Range("J64").Value = 4666.6665
$ws.Range("K64").Value = 3860
$ws.Range("L64").Value = 4666.6665
$ws.Range("M64").Value = -3612
$ws.Range("N64").Value = -5162.6665
$ws.Range("H67").Value = 4378.5713
$ws.Range("I67").Value = 3860
$ws.Range("J67").Value = 4666.6665
$ws.Range("K67").Value = 3860
$ws.Range("L67").Value = 4666.6665
$ws.Range("M67").Value = -3002
$ws.Range("N67").Value = -6382.6665
$ws.Range("H74").Value = 4319
$ws.Range("H77").Value = 4319
$ws.Range("H80").Value = 988883.4
$ws.Range("J80").Value = 813.25
$ws.Range("L80").Value = 2439.75
$ws.Range("N80").Value = -4435.75
$ws.Range("H83").Value = 988883.4
$ws.Range("J83").Value = 813.25
$ws.Range("L83").Value = 7319.25
$ws.Range("N83").Value = -17303.25
$ws.Range("H86").Value = 9117897
$ws.Range("J86").Value = 16712580
$ws.Range("L86").Value = 16712580
$ws.Range("N86").Value = -16714826
$ws.Range("H89").Value = 9117897
$ws.Range("J89").Value = 16712580
$ws.Range("L89").Value = 83562900
$ws.Range("N89").Value = -83574132
$ws.Range("H103").Value = 689.80646
$ws.Range("I103").Value = 413.8
$ws.Range("J103").Value = 821.2381
$ws.Range("K103").Value = 1241.4
$ws.Range("L103").Value = 2463.7143
$ws.Range("M103").Value = -655.4000000000001
$ws.Range("N103").Value = -3635.7143
$ws.Range("H111").Value = 23183.455
$ws.Range("I111").Value = 14494.375
$ws.Range("J111").Value = 46354.332
$ws.Range("K111").Value = 43483.125
$ws.Range("L111").Value = 139062.996
$ws.Range("M111").Value = -40416.125
$ws.Range("N111").Value = -145196.996
$ws.Range("H135").Value = 1305.0625
$ws.Range("I135").Value = 1356.7333
$ws.Range("J135").Value = 530
$ws.Range("K135").Value = 12210.5997
$ws.Range("L135").Value = 4770
$ws.Range("M135").Value = -9675.599700000001
$ws.Range("N135").Value = -9840
$ws.Range("H137").Value = 2394.02
$ws.Range("I137").Value = 1809.8182
$ws.Range("J137").Value = 2466.2246
$ws.Range("K137").Value = 5429.4546
$ws.Range("L137").Value = 7398.6738
$ws.Range("M137").Value = -2879.4546
$ws.Range("N137").Value = -12498.6738
$ws.Range("H138").Value = 2854
$ws.Range("I138").Value = 1840.7693
$ws.Range("K138").Value = 5522.3079
$ws.Range("M138").Value = -382.3078999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 63165400
$ws.Range("J63").Value = 28579642
$ws.Range("L63").Value = 28579642
$ws.Range("N63").Value = -28581014
$ws.Range("H66").Value = 63165400
$ws.Range("J66").Value = 28579642
$ws.Range("L66").Value = 142898210
$ws.Range("N66").Value = -142905074
$ws.Range("H110").Value = 52632732
$ws.Range("I110").Value = 58824640
$ws.Range("J110").Value = 1506.5
$ws.Range("K110").Value = 58824640
$ws.Range("L110").Value = 1506.5
$ws.Range("M110").Value = -58822595
$ws.Range("N110").Value = -5596.5
$ws.Range("H132").Value = 742539.0600000001
$ws.Range("I132").Value = 501341.7
$ws.Range("J132").Value = 1431674.4
$ws.Range("K132").Value = 1504025.1
$ws.Range("L132").Value = 4295023.199999999
$ws.Range("M132").Value = -1501495.1
$ws.Range("N132").Value = -4300083.199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1658.5333
$ws.Range("I86").Value = 1698.1666
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1698.1666
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -575.1666
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 1658.5333
$ws.Range("I89").Value = 1698.1666
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 8490.833000000001
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -2874.833000000001
$ws.Range("N89").Value = -18732
$ws.Range("H134").Value = 3593.2
$ws.Range("J134").Value = 3761.8
$ws.Range("L134").Value = 11285.4
$ws.Range("N134").Value = -16355.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2477.1633
$ws.Range("I31").Value = 732.6667
$ws.Range("J31").Value = 2869.675
$ws.Range("K31").Value = 732.6667
$ws.Range("L31").Value = 2869.675
$ws.Range("M31").Value = -437.6667
$ws.Range("N31").Value = -3459.675
$ws.Range("H34").Value = 2477.1633
$ws.Range("I34").Value = 732.6667
$ws.Range("J34").Value = 2869.675
$ws.Range("K34").Value = 732.6667
$ws.Range("L34").Value = 2869.675
$ws.Range("M34").Value = -530.6667
$ws.Range("N34").Value = -3273.675
$ws.Range("H62").Value = 8438.846
$ws.Range("I62").Value = 8308.75
$ws.Range("K62").Value = 8308.75
$ws.Range("M62").Value = -7684.75
$ws.Range("H65").Value = 8438.846
$ws.Range("I65").Value = 8308.75
$ws.Range("K65").Value = 41543.75
$ws.Range("M65").Value = -38423.75
$ws.Range("H109").Value = 45460
$ws.Range("J109").Value = 47946.668
$ws.Range("L109").Value = 47946.668
$ws.Range("N109").Value = -50026.668
$ws.Range("H132").Value = 2398.9167
$ws.Range("I132").Value = 2398.9167
$ws.Range("K132").Value = 7196.750100000001
$ws.Range("M132").Value = -4666.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1726.9166
$ws.Range("I68").Value = 1107.1428
$ws.Range("J68").Value = 1982.1177
$ws.Range("K68").Value = 3321.4284
$ws.Range("L68").Value = 5946.3531
$ws.Range("M68").Value = -2510.4284
$ws.Range("N68").Value = -7568.3531
$ws.Range("H71").Value = 1726.9166
$ws.Range("I71").Value = 1107.1428
$ws.Range("J71").Value = 1982.1177
$ws.Range("K71").Value = 9964.2852
$ws.Range("L71").Value = 17839.0593
$ws.Range("M71").Value = -5908.2852
$ws.Range("N71").Value = -25951.0593
$ws.Range("H107").Value = 1163.4
$ws.Range("I107").Value = 1163.4
$ws.Range("K107").Value = 3490.2
$ws.Range("M107").Value = -1570.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1139.8889
$ws.Range("I2").Value = 2531.5
$ws.Range("K2").Value = 2531.5
$ws.Range("M2").Value = -2418.5
$ws.Range("H62").Value = 57363.4
$ws.Range("J62").Value = 100000
$ws.Range("L62").Value = 100000
$ws.Range("N62").Value = -101372
$ws.Range("H65").Value = 57363.4
$ws.Range("J65").Value = 100000
$ws.Range("L65").Value = 300000
$ws.Range("N65").Value = -306864
$ws.Range("H126").Value = 8522.478999999999
$ws.Range("I126").Value = 11709
$ws.Range("J126").Value = 4380
$ws.Range("K126").Value = 35127
$ws.Range("L126").Value = 13140
$ws.Range("M126").Value = -32657
$ws.Range("N126").Value = -18080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3055.2942
$ws.Range("I40").Value = 2495.4167
$ws.Range("K40").Value = 2495.4167
$ws.Range("M40").Value = -2359.4167
$ws.Range("H132").Value = 5764.185
$ws.Range("I132").Value = 3253.5625
$ws.Range("J132").Value = 9416
$ws.Range("K132").Value = 9760.6875
$ws.Range("L132").Value = 28248
$ws.Range("M132").Value = -7230.6875
$ws.Range("N132").Value = -33308

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 416
$ws.Range("J100").Value = 351
$ws.Range("L100").Value = 702
$ws.Range("N100").Value = -1784
$ws.Range("H113").Value = 5033.923
$ws.Range("J113").Value = 7398.3335
$ws.Range("L113").Value = 22195.0005
$ws.Range("N113").Value = -26535.0005
$ws.Range("H132").Value = 574860.7
$ws.Range("I132").Value = 669854.1
$ws.Range("J132").Value = 4899.8
$ws.Range("K132").Value = 2009562.3
$ws.Range("L132").Value = 14699.4
$ws.Range("M132").Value = -2007032.3
$ws.Range("N132").Value = -19759.4
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("H136").Value = 4482.324
$ws.Range("I136").Value = 3715.3794
$ws.Range("J136").Value = 7262.5
$ws.Range("K136").Value = 11146.1382
$ws.Range("L136").Value = 21787.5
$ws.Range("M136").Value = -8596.138199999999
$ws.Range("N136").Value = -26887.5

# Remove M135 on WVR entirely (cell no longer present after update)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M135").ClearContents()
